$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Question text (used to sit in column A) moves to column B; row 3's question
# text is also updated at the same time ("is it electronic").
$questions = @(
    "organic or not",
    "do you use it",
    "is it electronic",
    "it is a place",
    "is it on earth",
    "is it a mammal",
    "can it fly",
    "does it walk on two legs"
)

# Item names (used to sit in column F) move to column G.
$items = @(
    "plate",
    "phone",
    "window",
    "city",
    "frog",
    "bird",
    "human",
    "star",
    "dog"
)

for ($r = 1; $r -le 8; $r++) {
    $ws.Cells.Item($r, 1).Value = $r - 1              # column A: numeric index 0..7
    $ws.Cells.Item($r, 2).Value = $questions[$r - 1]  # column B: question text
    $ws.Cells.Item($r, 6).ClearContents()             # column F: vacated
    $ws.Cells.Item($r, 7).Value = $items[$r - 1]      # column G: item text
}
$ws.Cells.Item(9, 6).ClearContents()
$ws.Cells.Item(9, 7).Value = $items[8]

# Tab color: alpha channel normalized from 00FFFFFF to FFFFFFFF (opaque white).
$ws.Tab.Color = 16777215

# Active selection moves to F8.
$ws.Range("F8").Select()
